$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Jrue Holiday -> Tyrese Maxey (PG,SG / Philadelphia 76ers)
$ws.Range("A3").Value = "Tyrese Maxey"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Philadelphia 76ers"

# Row 14: Tyrese Haliburton -> Jrue Holiday (PG,SG / Boston Celtics)
$ws.Range("A14").Value = "Jrue Holiday"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Boston Celtics"

# Row 15: Keegan Murray -> Tyrese Haliburton (PG,SG / Indiana Pacers)
$ws.Range("A15").Value = "Tyrese Haliburton"
$ws.Range("B15").Value = "PG,SG"
$ws.Range("C15").Value = "Indiana Pacers"

# Row 16: Tyrese Maxey -> Keegan Murray (SF,PF / Sacramento Kings)
$ws.Range("A16").Value = "Keegan Murray"
$ws.Range("B16").Value = "SF,PF"
$ws.Range("C16").Value = "Sacramento Kings"

# Row 18: Dereck Lively II -> Jalen Johnson (SF,PF / Atlanta Hawks)
$ws.Range("A18").Value = "Jalen Johnson"
$ws.Range("B18").Value = "SF,PF"
$ws.Range("C18").Value = "Atlanta Hawks"

# Row 19: Jalen Johnson -> Dereck Lively II (C / Dallas Mavericks)
$ws.Range("A19").Value = "Dereck Lively II"
$ws.Range("B19").Value = "C"
$ws.Range("C19").Value = "Dallas Mavericks"
